# Update "想去人数" (F column) values on both the "展览" and "全部类型"
# sheets, which carry duplicate data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 477
    $ws.Range("F3").Value = 3330
    $ws.Range("F4").Value = 86
    $ws.Range("F5").Value = 659
}
